$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1 (Professional Profile paragraph):
#   "currently completing a Master" -> "completed a Master"
#   Final XML keeps the text split as: "...comple" | "ted" | " a Master..."
# -----------------------------------------------------------------

# Step 1: remove "currently " (plain edit, no tracked-revision noise) so the
# surrounding text folds back into a single run.
$rng = $d.Content
$rng.Find.Execute("currently completing", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null
$start = $rng.Start
$d.Range($start, $start + 10).Delete()

# Step 2: locate "completing a Master" (now "currently " is gone) so we can
# change only the "ting" tail to "ted", isolating that tiny edit in its own
# run while "comple" stays merged with the preceding sentence.
$rng2 = $d.Content
$rng2.Find.Execute("completing a Master", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0) | Out-Null
$cStart = $rng2.Start
$tingRange = $d.Range($cStart + 6, $cStart + 10)

$d.TrackRevisions = $true
$tingRange.Find.Execute("ting", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ted", 2) | Out-Null
$d.AcceptAllRevisions()
$d.TrackRevisions = $false

# -----------------------------------------------------------------
# Change 2 (Relevant Courses cell): no textual change, but the run that held
#   "for Data and Analysis, Fundamentals of Data Engineering, ..."
# gets split right before "Engineering" into two runs.
# -----------------------------------------------------------------

$d.TrackRevisions = $true
$d.Content.Find.Execute(
    "Engineering, Applied Machine Learning, Statistical Methods for Discrete Response, Time Series, and Panel Data",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Engineering, Applied Machine Learning, Statistical Methods for Discrete Response, Time Series, and Panel Data",
    2) | Out-Null
$d.AcceptAllRevisions()
$d.TrackRevisions = $false
